$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("Settings")
$assets = $wb.Worksheets.Item("Assets")

# --- Settings sheet: update existing rows -----------------------------
# Row 2: Orchestrator queue name / value
$settings.Cells.Item(2,1).Value = "OrchestratorQueueName"
$settings.Cells.Item(2,2).Value = "P004_SP002_090_NHC_OLC_Dispatcher_Queue"
$settings.Cells.Item(2,3).Value = "Orchestrator queue Name. The value must match with the queue name defined on Orchestrator."

# Row 3: Orchestrator queue folder
$settings.Cells.Item(3,1).Value = "OrchestratorQueueFolder"
$settings.Cells.Item(3,2).Value = "DEV"
$settings.Cells.Item(3,3).Value = "Folder name. The value must match a folder defined in Orchestrator and queue specified as OrchestratorQueueName should be created in this folder. For classic folders leave the value field empty."
$settings.Rows.Item(3).RowHeight = 45

# Row 5: business process name
$settings.Cells.Item(5,1).Value = "logF_BusinessProcessName"
$settings.Cells.Item(5,2).Value = "P004_SP002_090_NHC_OLC_Dispatcher"
$settings.Cells.Item(5,3).Value = "Logging field which allows grouping of log data of two or more subprocesses under the same business process name"
$settings.Rows.Item(5).RowHeight = 30

# --- Settings sheet: new rows for SharePoint dispatcher config --------
# Row 8: SharePoint URL (with hyperlink)
$settings.Cells.Item(8,1).Value = "SharePointURL"
$settings.Cells.Item(8,2).Value = "https://officemgmtentserv.sharepoint.com/sites/NewHireCommunication/Lists/New%20Hire%20Communication%20%20Employee%20Details/OfferLetterTestView.aspx"
$settings.Cells.Item(8,2).WrapText = $true
$settings.Hyperlinks.Add($settings.Cells.Item(8,2), "https://officemgmtentserv.sharepoint.com/sites/NewHireCommunication/Lists/New%20Hire%20Communication%20%20Employee%20Details/OfferLetterTestView.aspx")
$settings.Cells.Item(8,3).Value = "URL of SharepointData Input"

# Row 9: Offer letter input excel
$settings.Cells.Item(9,1).Value = "OfferLetterInputExcel"
$settings.Cells.Item(9,2).Value = "C:\Users\55649C\Documents\UiPath\P004_SP002_090_NewHireCommunication_OfferLetterCreation_Dispatcher\Data\Input\OfferLetterInput.xlsx"
$settings.Cells.Item(9,3).Value = "Dummy Data for Sharepoint"

# Row 10: Offer letter table name
$settings.Cells.Item(10,1).Value = "OfferLetterTableName"
$settings.Cells.Item(10,2).Value = "OfferLetterList"

# Row 12: Name extractor regex
$settings.Cells.Item(12,1).Value = "NameExtractor"
$settings.Cells.Item(12,2).Value = "applicant.|(First)"

# Row 13: Email extractor regex
$settings.Cells.Item(13,1).Value = "EmailExtractor"
$settings.Cells.Item(13,2).Value = "Birth|(Email"

# --- Assets sheet: page orientation ------------------------------------
$assets.PageSetup.Orientation = 1

# --- Sheet selection / activation --------------------------------------
$settings.Activate()
$settings.Range("C16").Select()
